# Update the NATMI ligand-receptor pair sheet (Thbs1-Cd36) with refreshed TPM-derived
# expression / specificity / weight values (rows 2-17, columns G:J and M:T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("M2").Value = 509.3923236666667
$ws.Range("N2").Value = 1528.176971
$ws.Range("O2").Value = 0.831019558191033
$ws.Range("P2").Value = 0.8310195581910331
$ws.Range("Q2").Value = 10791.98559940456
$ws.Range("R2").Value = 97127.87039464104
$ws.Range("S2").Value = 0.07284004100175223
$ws.Range("T2").Value = 0.07284004100175225

$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("N3").Value = 3.779073
$ws.Range("O3").Value = 0.002055052284145212
$ws.Range("P3").Value = 0.002055052284145212
$ws.Range("Q3").Value = 26.687813106103
$ws.Range("R3").Value = 240.190317954927
$ws.Range("S3").Value = 0.0001801282426658259
$ws.Range("T3").Value = 0.000180128242665826

$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 22.92703233333333
$ws.Range("N4").Value = 68.781097
$ws.Range("O4").Value = 0.03740302198339736
$ws.Range("P4").Value = 0.03740302198339737
$ws.Range("Q4").Value = 485.7320993716559
$ws.Range("R4").Value = 4371.588894344904
$ws.Range("S4").Value = 0.00327842783964155
$ws.Range("T4").Value = 0.003278427839641551

$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("M5").Value = 79.39367866666667
$ws.Range("N5").Value = 238.181036
$ws.Range("O5").Value = 0.1295223675414243
$ws.Range("P5").Value = 0.1295223675414243
$ws.Range("Q5").Value = 1682.034449767441
$ws.Range("R5").Value = 15138.31004790697
$ws.Range("S5").Value = 0.01135281891908567
$ws.Range("T5").Value = 0.01135281891908567

$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("M6").Value = 509.3923236666667
$ws.Range("N6").Value = 1528.176971
$ws.Range("O6").Value = 0.831019558191033
$ws.Range("P6").Value = 0.8310195581910331
$ws.Range("Q6").Value = 62054.6274168235
$ws.Range("R6").Value = 558491.6467514115
$ws.Range("S6").Value = 0.41883502936098
$ws.Range("T6").Value = 0.41883502936098

$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("N7").Value = 3.779073
$ws.Range("O7").Value = 0.002055052284145212
$ws.Range("P7").Value = 0.002055052284145212
$ws.Range("S7").Value = 0.001035749249562724
$ws.Range("T7").Value = 0.001035749249562724

$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 22.92703233333333
$ws.Range("N8").Value = 68.781097
$ws.Range("O8").Value = 0.03740302198339736
$ws.Range("P8").Value = 0.03740302198339737
$ws.Range("Q8").Value = 2792.991537401852
$ws.Range("R8").Value = 25136.92383661667
$ws.Range("S8").Value = 0.01885117583117631
$ws.Range("T8").Value = 0.01885117583117631

$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("M9").Value = 79.39367866666667
$ws.Range("N9").Value = 238.181036
$ws.Range("O9").Value = 0.1295223675414243
$ws.Range("P9").Value = 0.1295223675414243
$ws.Range("Q9").Value = 9671.808780799265
$ws.Range("R9").Value = 87046.27902719338
$ws.Range("S9").Value = 0.06527945591341373
$ws.Range("T9").Value = 0.06527945591341373

$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("M10").Value = 509.3923236666667
$ws.Range("N10").Value = 1528.176971
$ws.Range("O10").Value = 0.831019558191033
$ws.Range("P10").Value = 0.8310195581910331
$ws.Range("Q10").Value = 18953.05306586601
$ws.Range("R10").Value = 170577.4775927941
$ws.Range("S10").Value = 0.1279228136203441
$ws.Range("T10").Value = 0.1279228136203441

$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("N11").Value = 3.779073
$ws.Range("O11").Value = 0.002055052284145212
$ws.Range("P11").Value = 0.002055052284145212
$ws.Range("Q11").Value = 46.86955272065899
$ws.Range("R11").Value = 421.8259744859309
$ws.Range("S11").Value = 0.0003163440231142407
$ws.Range("T11").Value = 0.0003163440231142409

$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 22.92703233333333
$ws.Range("N12").Value = 68.781097
$ws.Range("O12").Value = 0.03740302198339736
$ws.Range("P12").Value = 0.03740302198339737
$ws.Range("Q12").Value = 853.050272388562
$ws.Range("R12").Value = 7677.452451497059
$ws.Range("S12").Value = 0.005757625994308878
$ws.Range("T12").Value = 0.005757625994308881

$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("M13").Value = 79.39367866666667
$ws.Range("N13").Value = 238.181036
$ws.Range("O13").Value = 0.1295223675414243
$ws.Range("P13").Value = 0.1295223675414243
$ws.Range("Q13").Value = 2954.015078264743
$ws.Range("R13").Value = 26586.13570438269
$ws.Range("S13").Value = 0.01993799727016594
$ws.Range("T13").Value = 0.01993799727016594

$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("M14").Value = 509.3923236666667
$ws.Range("N14").Value = 1528.176971
$ws.Range("O14").Value = 0.831019558191033
$ws.Range("P14").Value = 0.8310195581910331
$ws.Range("Q14").Value = 31324.25012500177
$ws.Range("R14").Value = 281918.2511250159
$ws.Range("S14").Value = 0.2114216742079567
$ws.Range("T14").Value = 0.2114216742079567

$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("N15").Value = 3.779073
$ws.Range("O15").Value = 0.002055052284145212
$ws.Range("P15").Value = 0.002055052284145212
$ws.Range("Q15").Value = 77.46264348897898
$ws.Range("R15").Value = 697.1637914008109
$ws.Range("S15").Value = 0.0005228307688024213
$ws.Range("T15").Value = 0.0005228307688024214

$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 22.92703233333333
$ws.Range("N16").Value = 68.781097
$ws.Range("O16").Value = 0.03740302198339736
$ws.Range("P16").Value = 0.03740302198339737
$ws.Range("Q16").Value = 1409.860459348597
$ws.Range("R16").Value = 12688.74413413738
$ws.Range("S16").Value = 0.009515792318270622
$ws.Range("T16").Value = 0.009515792318270624

$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("M17").Value = 79.39367866666667
$ws.Range("N17").Value = 238.181036
$ws.Range("O17").Value = 0.1295223675414243
$ws.Range("P17").Value = 0.1295223675414243
$ws.Range("Q17").Value = 4882.184778516761
$ws.Range("R17").Value = 43939.66300665084
$ws.Range("S17").Value = 0.03295209543875897
$ws.Range("T17").Value = 0.03295209543875897
